$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.732.78"
$ws.Range("E2").Value = "  +1.30%  "
$ws.Range("D3").Value = "1.725.99"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("D4").Value = "0.9979"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "240.45"
$ws.Range("E5").Value = "  -0.89%  "
$ws.Range("D6").Value = "0.9985"
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").Value = "0.4836"
$ws.Range("E7").Value = "  -0.92%  "
$ws.Range("D8").Value = "0.2585"
$ws.Range("E8").Value = "  -0.25%  "
$ws.Range("D9").Value = "0.06184"
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("D10").Value = "1.727.37"
$ws.Range("E10").Value = "  +0.15%  "
$ws.Range("D11").Value = "15.88"
$ws.Range("E11").Value = "  +2.25%  "
$ws.Range("D12").Value = "0.06874"
$ws.Range("E12").Value = "  -1.53%  "
$ws.Range("D13").Value = "0.6041"
$ws.Range("E13").Value = "  +1.05%  "
$ws.Range("D14").Value = "4.466"
$ws.Range("E14").Value = "  -1.33%  "
$ws.Range("D15").Value = "76.96"
$ws.Range("E15").Value = "  -0.29%  "
$ws.Range("D16").Value = "0.9984"
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("D17").Value = "26.556.44"
$ws.Range("E17").Value = "  +0.62%  "
$ws.Range("D18").Value = "0.9981"
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("D19").Value = "0.000007164"
$ws.Range("E19").Value = "  -0.55%  "
$ws.Range("D20").Value = "11.36"
$ws.Range("E20").Value = "  +0.40%  "
$ws.Range("D21").Value = "1.945.72"
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").Value = "4.417"
$ws.Range("E22").Value = "  -0.65%  "
$ws.Range("D23").Value = "8.543"
$ws.Range("E23").Value = "  +0.56%  "
$ws.Range("D24").Value = "5.055"
$ws.Range("E24").Value = "  -0.92%  "
$ws.Range("D25").Value = "140.17"
$ws.Range("E25").Value = "  +1.47%  "
$ws.Range("D26").Value = "15.23"
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("D27").Value = "1.778"
$ws.Range("E27").Value = "  +2.98%  "
$ws.Range("D28").Value = "106.36"
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("D29").Value = "1.368"
$ws.Range("E29").Value = "  -2.27%  "
$ws.Range("D30").Value = "4.015"
$ws.Range("E30").Value = "  +2.79%  "
$ws.Range("D31").Value = "0.07933"
$ws.Range("E31").Value = "  -0.99%  "
$ws.Range("D32").Value = "3.670"
$ws.Range("E32").Value = "  +0.39%  "
$ws.Range("D33").Value = "0.04515"
$ws.Range("E33").Value = "  +0.35%  "
$ws.Range("D34").Value = "2.596"
$ws.Range("E34").Value = "  -0.28%  "
$ws.Range("D35").Value = "1.001"
$ws.Range("E35").Value = "  +0.17%  "
$ws.Range("D36").Value = "0.6187"
$ws.Range("E36").Value = "  -0.75%  "
$ws.Range("D37").Value = "0.9355"
$ws.Range("E37").Value = "  +0.56%  "
$ws.Range("D38").Value = "2.004"
$ws.Range("E38").Value = "  +2.22%  "
$ws.Range("D39").Value = "2.451"
$ws.Range("E39").Value = "  +3.03%  "
$ws.Range("D40").Value = "0.9978"
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("D41").Value = "0.01498"
$ws.Range("E41").Value = "  +1.64%  "
$ws.Range("D42").Value = "5.608"
$ws.Range("E42").Value = "  +2.56%  "
$ws.Range("D43").Value = "99.79"
$ws.Range("E43").Value = "  -0.78%  "
$ws.Range("D44").Value = "0.3833"
$ws.Range("E44").Value = "  -0.27%  "
$ws.Range("D45").Value = "6.797"
$ws.Range("E45").Value = "  -1.84%  "
$ws.Range("D46").Value = "0.1154"
$ws.Range("E46").Value = "  -0.93%  "
$ws.Range("D47").Value = "0.05358"
$ws.Range("E47").Value = "  -0.19%  "
$ws.Range("D48").Value = "7.950"
$ws.Range("E48").Value = "  +3.73%  "
$ws.Range("D49").Value = "30.10"
$ws.Range("E49").Value = "  -0.21%  "
$ws.Range("D50").Value = "1.244"
$ws.Range("E50").Value = "  +1.46%  "
$ws.Range("D51").Value = "51.42"
$ws.Range("E51").Value = "  +0.78%  "
